$wb = $excel.ActiveWorkbook

$wsProduct = $wb.Worksheets.Item("ProductDetails")

# Update row 2: Purse -> pens, quantity 5 -> 2 (Sort stays "Price: Low to High")
$wsProduct.Range("A2").Value = "pens"
$wsProduct.Range("B2").Value = 2

# Delete row 3 entirely (Mobiles | 4 | Newest Arrivals)
$wsProduct.Rows.Item(3).Delete()
